$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List_ID")

$ws.Range("C16").Value = "14/08/1997"
$ws.Range("D16").Value = "xaprkkcwssjkbsl@gmail.com"
$ws.Range("E16").Value = "neddhSAIKM5"
$ws.Range("F16").Value = "pass"

$ws.Range("C17").Value = "17/01/1984"
$ws.Range("D17").Value = "bjveamxemxbginv@gmail.com"
$ws.Range("E17").Value = "mvmugZZFRW5"
$ws.Range("F17").Value = "pass"
